$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1829652996845426
$ws.Range("C2").Value = 0.5678233438485805
$ws.Range("J2").Value = 0.01577287066246057
$ws.Range("P2").Value = 0.138801261829653
$ws.Range("S2").Value = 0.0946372239747634
$ws.Range("B3").Value = 0.01621621621621622
$ws.Range("C3").Value = 0.02702702702702703
$ws.Range("J3").Value = 0.04324324324324325
$ws.Range("P3").Value = 0.772972972972973
$ws.Range("S3").Value = 0.1405405405405405
$ws.Range("J4").Value = 0.1395348837209302
$ws.Range("P4").Value = 0.6511627906976745
$ws.Range("S4").Value = 0.2093023255813954
$ws.Range("B6").Value = 0.08225108225108226
$ws.Range("F6").Value = 0.06926406926406926
$ws.Range("J6").Value = 0.2121212121212121
$ws.Range("O6").Value = 0.01298701298701299
$ws.Range("Q6").Value = 0.1688311688311688
$ws.Range("R6").Value = 0.06493506493506493
$ws.Range("S6").Value = 0.3896103896103896
$ws.Range("B7").Value = 0.1234567901234568
$ws.Range("D7").Value = 0.01851851851851852
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.154320987654321
$ws.Range("O7").Value = 0.04320987654320987
$ws.Range("Q7").Value = 0.191358024691358
$ws.Range("R7").Value = 0.07407407407407407
$ws.Range("S7").Value = 0.3395061728395062
$ws.Range("B8").Value = 0.1063829787234043
$ws.Range("D8").Value = 0.02659574468085106
$ws.Range("F8").Value = 0.07180851063829788
$ws.Range("J8").Value = 0.1223404255319149
$ws.Range("O8").Value = 0.02393617021276596
$ws.Range("Q8").Value = 0.1781914893617021
$ws.Range("R8").Value = 0.1037234042553191
$ws.Range("S8").Value = 0.3670212765957447
$ws.Range("B9").Value = 0.08379888268156424
$ws.Range("D9").Value = 0.0223463687150838
$ws.Range("F9").Value = 0.1005586592178771
$ws.Range("J9").Value = 0.08379888268156424
$ws.Range("O9").Value = 0.0446927374301676
$ws.Range("Q9").Value = 0.1620111731843575
$ws.Range("R9").Value = 0.1005586592178771
$ws.Range("S9").Value = 0.4022346368715084
$ws.Range("B10").Value = 0.1282258064516129
$ws.Range("D10").Value = 0.0217741935483871
$ws.Range("E10").Value = 0.0008064516129032258
$ws.Range("F10").Value = 0.08064516129032258
$ws.Range("J10").Value = 0.1145161290322581
$ws.Range("O10").Value = 0.02016129032258064
$ws.Range("Q10").Value = 0.2169354838709678
$ws.Range("R10").Value = 0.07096774193548387
$ws.Range("S10").Value = 0.3459677419354839
$ws.Range("G11").Value = 0.1286764705882353
$ws.Range("J11").Value = 0.1102941176470588
$ws.Range("K11").Value = 0.2132352941176471
$ws.Range("L11").Value = 0.5367647058823529
$ws.Range("S11").Value = 0.01102941176470588
$ws.Range("G12").Value = 0.7152317880794702
$ws.Range("J12").Value = 0.2251655629139073
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.02649006622516556
$ws.Range("S12").Value = 0.02649006622516556
$ws.Range("F15").Value = 0.02
$ws.Range("H15").Value = 0.128
$ws.Range("I15").Value = 0.048
$ws.Range("J15").Value = 0.38
$ws.Range("K15").Value = 0.068
$ws.Range("M15").Value = 0.012
$ws.Range("O15").Value = 0.076
$ws.Range("S15").Value = 0.268
$ws.Range("F16").Value = 0.01428571428571429
$ws.Range("H16").Value = 0.119047619047619
$ws.Range("I16").Value = 0.05714285714285714
$ws.Range("J16").Value = 0.4619047619047619
$ws.Range("K16").Value = 0.1
$ws.Range("M16").Value = 0.02380952380952381
$ws.Range("O16").Value = 0.07142857142857142
$ws.Range("S16").Value = 0.1523809523809524
$ws.Range("F17").Value = 0.0207852193995381
$ws.Range("H17").Value = 0.1639722863741339
$ws.Range("I17").Value = 0.09006928406466513
$ws.Range("J17").Value = 0.4087759815242494
$ws.Range("K17").Value = 0.115473441108545
$ws.Range("M17").Value = 0.01847575057736721
$ws.Range("N17").Value = 0.002309468822170901
$ws.Range("O17").Value = 0.07852193995381063
$ws.Range("S17").Value = 0.1016166281755196
$ws.Range("F18").Value = 0.01775147928994083
$ws.Range("H18").Value = 0.2071005917159763
$ws.Range("I18").Value = 0.05325443786982249
$ws.Range("J18").Value = 0.4792899408284024
$ws.Range("K18").Value = 0.05917159763313609
$ws.Range("M18").Value = 0.01183431952662722
$ws.Range("O18").Value = 0.07100591715976332
$ws.Range("S18").Value = 0.1005917159763314
$ws.Range("F19").Value = 0.01422222222222222
$ws.Range("H19").Value = 0.1884444444444444
$ws.Range("I19").Value = 0.09333333333333334
$ws.Range("J19").Value = 0.3946666666666667
$ws.Range("K19").Value = 0.09866666666666667
$ws.Range("M19").Value = 0.01422222222222222
$ws.Range("O19").Value = 0.08177777777777778
$ws.Range("S19").Value = 0.1146666666666667
